# Apply the cryptos-list refresh described in the commit:
# "Updated cryptos list on Thu Oct 26 23:34:15 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some new Price values are plain decimals (e.g. "223.34") that Excel's
# COM layer would otherwise auto-coerce to a Number. Force those specific
# cells to Text format first so the assigned string is preserved verbatim,
# exactly like the rest of the (already-text) Price column.
$textFormatCells = @(
    "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D14", "D19", "D21",
    "D23", "D25", "D26", "D30", "D39", "D40", "D41", "D42", "D44", "D45",
    "D49", "D50"
)
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Write the updated coin figures.
$ws.Range("D2").Value = "34.205.70"
$ws.Range("E2").Value = "  -0.89%  "
$ws.Range("D3").Value = "1.806.70"
$ws.Range("E3").Value = "  +1.07%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "223.34"
$ws.Range("E5").Value = "  +0.43%  "
$ws.Range("D6").Value = "0.553"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "32.98"
$ws.Range("E8").Value = "  +2.07%  "
$ws.Range("D9").Value = "0.287"
$ws.Range("E9").Value = "  +2.51%  "
$ws.Range("D10").Value = "0.0717"
$ws.Range("D11").Value = "0.0929"
$ws.Range("E11").Value = "  -0.66%  "
$ws.Range("D12").Value = "2.067.87"
$ws.Range("E12").Value = "  +1.15%  "
$ws.Range("D13").Value = "1.819.39"
$ws.Range("E13").Value = "  +1.74%  "
$ws.Range("D14").Value = "11.02"
$ws.Range("E14").Value = "  +0.34%  "
$ws.Range("E15").Value = "  +0.29%  "
$ws.Range("D16").Value = "34.240.41"
$ws.Range("E16").Value = "  -0.78%  "
$ws.Range("E17").Value = "  -0.92%  "
$ws.Range("E18").Value = "  +0.37%  "
$ws.Range("D19").Value = "247.63"
$ws.Range("E19").Value = "  -2.46%  "
$ws.Range("D20").Value = "0.0₃0789"
$ws.Range("E20").Value = "  +1.11%  "
$ws.Range("D21").Value = "11.12"
$ws.Range("E21").Value = "  +6.10%  "
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").Value = "4.11"
$ws.Range("E23").Value = "  -1.24%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").Value = "159.77"
$ws.Range("E25").Value = "  -0.38%  "
$ws.Range("D26").Value = "16.60"
$ws.Range("E26").Value = "  +1.33%  "
$ws.Range("E27").Value = "  +1.09%  "
$ws.Range("E28").Value = "  -0.96%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").Value = "0.0529"
$ws.Range("E30").Value = "  +2.60%  "
$ws.Range("E31").Value = "  -0.34%  "
$ws.Range("E32").Value = "  +1.57%  "
$ws.Range("E33").Value = "  -0.45%  "
$ws.Range("E34").Value = "  -0.73%  "
$ws.Range("D35").Value = "1.421.08"
$ws.Range("E35").Value = "  -0.96%  "
$ws.Range("E36").Value = "  +2.40%  "
$ws.Range("E37").Value = "  +0.52%  "
$ws.Range("E38").Value = "  -1.06%  "
$ws.Range("D39").Value = "0.944"
$ws.Range("E39").Value = "  +3.01%  "
$ws.Range("D40").Value = "80.90"
$ws.Range("E40").Value = "  -4.91%  "
$ws.Range("B41").Value = "HuobiToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D41").Value = "2.36"
$ws.Range("E41").Value = "  +0.85%  "
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").Value = "2.73"
$ws.Range("E42").Value = "  -2.26%  "
$ws.Range("E43").Value = "  +4.16%  "
$ws.Range("D44").Value = "5.97"
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("D45").Value = "108.23"
$ws.Range("E45").Value = "  +4.17%  "
$ws.Range("E46").Value = "  +1.12%  "
$ws.Range("E47").Value = "  -0.96%  "
$ws.Range("D48").Value = "1.966.29"
$ws.Range("E48").Value = "  +1.22%  "
$ws.Range("D49").Value = "12.15"
$ws.Range("E49").Value = "  +0.83%  "
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("E51").Value = "  +2.92%  "
